$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '98.519.28'
$ws.Range("E2").Value = '  +1.10%  '

# Row 3
$ws.Range("D3").Value = '3.496.07'
$ws.Range("E3").Value = '  +4.26%  '

# Row 4
$ws.Range("E4").Value = '  -0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '253.60'
$ws.Range("E5").Value = '  +1.67%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '666.85'
$ws.Range("E6").Value = '  +1.89%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.49'
$ws.Range("E7").Value = '  +6.50%  '

# Row 8
$ws.Range("E8").Value = '  +2.21%  '

# Row 9
$ws.Range("E9").Value = '  +3.34%  '

# Row 10
$ws.Range("E10").Value = '  +0.00%  '

# Row 11
$ws.Range("D11").Value = '3.492.86'
$ws.Range("E11").Value = '  +4.22%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.74'
$ws.Range("E12").Value = '  +12.69%  '

# Row 13
$ws.Range("E13").Value = '  +1.00%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.28'
$ws.Range("E14").Value = '  +3.02%  '

# Row 15
$ws.Range("D15").Value = '98.201.92'
$ws.Range("E15").Value = '  +0.87%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000262'
$ws.Range("E16").Value = '  +2.63%  '

# Row 17
$ws.Range("D17").Value = '4.157.46'
$ws.Range("E17").Value = '  +4.34%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '9.01'
$ws.Range("E18").Value = '  +4.01%  '

# Row 19
$ws.Range("D19").Value = '3.508.97'
$ws.Range("E19").Value = '  +4.80%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.95'
$ws.Range("E20").Value = '  +12.38%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.83'
$ws.Range("E21").Value = '  +9.60%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.531'
$ws.Range("E22").Value = '  -4.78%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '527.69'
$ws.Range("E23").Value = '  +4.94%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.44'
$ws.Range("E24").Value = '  +2.57%  '

# Row 26
$ws.Range("E26").Value = '  +7.71%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '98.84'
$ws.Range("E27").Value = '  +1.97%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.80'
$ws.Range("E28").Value = '  +5.34%  '

# Row 29
$ws.Range("D29").Value = '3.687.29'
$ws.Range("E29").Value = '  +4.07%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.50'
$ws.Range("E30").Value = '  +12.74%  '

# Row 31
$ws.Range("E31").Value = '  +14.49%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.146'
$ws.Range("E32").Value = '  -0.97%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.997'
$ws.Range("E33").Value = '  +0.07%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.191'
$ws.Range("E34").Value = '  -0.20%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.611'
$ws.Range("E35").Value = '  +10.28%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.40'
$ws.Range("E36").Value = '  +9.47%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.49%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.53'
$ws.Range("E38").Value = '  +4.73%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.00'
$ws.Range("E39").Value = '  +3.36%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.157'
$ws.Range("E40").Value = '  +4.99%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '528.57'
$ws.Range("E41").Value = '  +1.67%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.928'
$ws.Range("E42").Value = '  +9.67%  '

# Row 43
$ws.Range("E43").Value = '  +0.09%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.78'
$ws.Range("E44").Value = '  +7.26%  '

# Row 45
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0438'
$ws.Range("E45").Value = '  +3.89%  '

# Row 46
$ws.Range("B46").Value = 'WhiteBITCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '24.42'
$ws.Range("E46").Value = '  -0.89%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.82'
$ws.Range("E47").Value = '  +3.40%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.76'
$ws.Range("E48").Value = '  -1.30%  '

# Row 49
$ws.Range("E49").Value = '  -1.02%  '

# Row 50
$ws.Range("E50").Value = '  +12.38%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '55.36'
$ws.Range("E51").Value = '  +3.88%  '
